$p = $ppt.ActivePresentation

# Slide 2 ("Objectives") holds the bullet list; its second shape is the
# body/content placeholder with the "Create custom objects..." bullet.
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)

# Fall back to a text-based lookup in case shape ordering ever differs.
if (-not ($shp.HasTextFrame -and $shp.TextFrame.HasText -and `
        $shp.TextFrame.TextRange.Text -like "Create custom objects if needed*")) {
    foreach ($candidate in $s.Shapes) {
        if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText -and `
                $candidate.TextFrame.TextRange.Text -like "Create custom objects if needed*") {
            $shp = $candidate
            break
        }
    }
}

$tr = $shp.TextFrame.TextRange
$para = $tr.Paragraphs(1, 1)

# Original single run read:
#   "Create custom objects if needed and custom fields for storing student's
#    information as well as quote information."
#
# Split it into three runs, inserting the word "create" so the sentence
# reads "...needed and create custom fields...":
#   1) "Create custom objects if needed "
#   2) "and create custom "
#   3) "fields for storing student's information as well as quote information."
$middle = $para.Characters(33, 11)
$middle.Text = "and create custom "
